# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on Sheet1 previously held a "Strike#" style count; the data was
# regenerated so that column G now holds the actual strikeout total ("K")
# for each start. Write the recalculated K values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 4
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 2
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
